# "charlotte added a slide"
#
# A brand-new "Title and Content" slide ("Charlotte made a change!") is
# inserted as the very first slide of the deck, pushing every existing
# slide back by one position (old slide N -> position N+1). The new
# slide's content placeholder is left empty, just like a layout that was
# never typed into.

$p = $ppt.ActivePresentation

# Insert the new slide at position 1 using the "Title and Content" layout
# (layout #2 in this deck's slide-layout list), matching the placeholder
# set (Title + Content Placeholder) seen on the new slide.
$newSlide = $p.Slides.Add(1, 2)

# Re-fetch slide 1 (the one we just inserted) and set its title text.
$title = $p.Slides.Item(1).Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Charlotte made "
$titleRange.InsertAfter("a change!")

# The content placeholder (shape 2) is left empty - nothing was typed
# into it, matching the source slide.
